$wb = $excel.ActiveWorkbook

# --- Salesforce sheet: insert a new "AccountName" column between
#     OpportunityName and OpportunityID, and clear the old OpportunityID
#     sample values (column shifts from C to D). ---
$wsSales = $wb.Worksheets.Item("Salesforce")
$wsSales.Activate()

$wsSales.Columns.Item(3).Insert()
$wsSales.Columns.Item(3).ColumnWidth = $wsSales.Columns.Item(2).ColumnWidth

$wsSales.Range("C1").Value = "AccountName"
$wsSales.Range("C2").Value = "Adams25 Inc"
$wsSales.Range("C3").Value = "Adams25 Inc"
$wsSales.Range("C4").Value = "Adams25 Inc"

$wsSales.Range("D2").Value = ""

$wsSales.Range("C14").Select()

# --- Ecommerce sheet: clear the sample ProductId / OrderID values on row 2 ---
$wsEcom = $wb.Worksheets.Item("Ecommerce")
$wsEcom.Activate()

$wsEcom.Range("C2").Value = ""
$wsEcom.Range("D2").Value = ""

$wsEcom.Range("B2").Select()
